{"js": "// Update the date title and every arithmetic answer cell in the single\n// table, matching the author's commit (positional replacement, since a\n// few \"before\" expressions repeat with different \"after\" values).\n\nconst body = context.document.body;\n\n// ---- 1. Title paragraph: \"2023-05-29 Monday\" -> \"2023-05-30 Tuesday\" ----\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\ntitlePara.load(\"text\");\nawait context.sync();\nif (titlePara.text.trim() === \"2023-05-29 Monday\") {\n  titlePara.insertText(\"2023-05-30 Tuesday\", \"Replace\");\n}\n\n// ---- 2. The 20x5 table of \"a+b=\" / \"a-b=\" problems ----\nconst newValues = [\n  [\"70+19=\", \"24+74=\", \"4+47=\", \"32+22=\", \"51-0=\"],\n  [\"83-32=\", \"13+42=\", \"26+20=\", \"85-28=\", \"10-8=\"],\n  [\"47-23=\", \"49+50=\", \"22+69=\", \"60+23=\", \"99-70=\"],\n  [\"0+34=\", \"58-41=\", \"15+20=\", \"26+2=\", \"62-28=\"],\n  [\"64-38=\", \"54+2=\", \"39-34=\", \"63+23=\", \"96-3=\"],\n  [\"19-13=\", \"71+20=\", \"25+55=\", \"68-17=\", \"48-26=\"],\n  [\"44+54=\", \"42+10=\", \"37+35=\", \"73-68=\", \"74-52=\"],\n  [\"14+80=\", \"30-13=\", \"6+70=\", \"44-37=\", \"19+36=\"],\n  [\"9+29=\", \"93-43=\", \"22+41=\", \"65-25=\", \"58-39=\"],\n  [\"12+43=\", \"99-69=\", \"26-3=\", \"83-50=\", \"35-16=\"],\n  [\"78-35=\", \"4+41=\", \"25+11=\", \"58+21=\", \"72+18=\"],\n  [\"84-3=\", \"11+38=\", \"86-72=\", \"65-31=\", \"75-3=\"],\n  [\"96-76=\", \"22+76=\", \"26+39=\", \"12+5=\", \"68+15=\"],\n  [\"34-17=\", \"80-77=\", \"0+92=\", \"94-34=\", \"3+18=\"],\n  [\"76-58=\", \"17+13=\", \"51-50=\", \"9+10=\", \"58+9=\"],\n  [\"20+63=\", \"34+5=\", \"1-0=\", \"6+75=\", \"64-47=\"],\n  [\"44-35=\", \"79-12=\", \"97-85=\", \"88-82=\", \"4+48=\"],\n  [\"20-16=\", \"37+51=\", \"5+71=\", \"94-43=\", \"56-14=\"],\n  [\"25+70=\", \"9+85=\", \"4+24=\", \"91-47=\", \"92-91=\"],\n  [\"69-11=\", \"47-4=\", \"42+20=\", \"42+25=\", \"34+22=\"]\n];\n\nconst table = body.tables.getFirst();\ntable.load(\"rowCount\");\nawait context.sync();\n\nfor (let r = 0; r < newValues.length; r++) {\n  for (let c = 0; c < newValues[r].length; c++) {\n    table.getCell(r, c).value = newValues[r][c];\n  }\n}\nawait context.sync();\n", "ps1": "# Update the date title and every arithmetic answer cell in the single\n# table, matching the author's commit (positional replacement, since a\n# few \"before\" expressions repeat with different \"after\" values).\n\n$d = $word.ActiveDocument\n\n# ---- 1. Title paragraph: \"2023-05-29 Monday\" -> \"2023-05-30 Tuesday\" ----\n$titlePara = $d.Paragraphs(1)\nif ($titlePara.Range.Text.Trim() -eq \"2023-05-29 Monday\") {\n    $titlePara.Range.Text = \"2023-05-30 Tuesday\"\n}\n\n# ---- 2. The 20x5 table of \"a+b=\" / \"a-b=\" problems ----\n$newValues = @(\n    @(\"70+19=\", \"24+74=\", \"4+47=\", \"32+22=\", \"51-0=\"),\n    @(\"83-32=\", \"13+42=\", \"26+20=\", \"85-28=\", \"10-8=\"),\n    @(\"47-23=\", \"49+50=\", \"22+69=\", \"60+23=\", \"99-70=\"),\n    @(\"0+34=\", \"58-41=\", \"15+20=\", \"26+2=\", \"62-28=\"),\n    @(\"64-38=\", \"54+2=\", \"39-34=\", \"63+23=\", \"96-3=\"),\n    @(\"19-13=\", \"71+20=\", \"25+55=\", \"68-17=\", \"48-26=\"),\n    @(\"44+54=\", \"42+10=\", \"37+35=\", \"73-68=\", \"74-52=\"),\n    @(\"14+80=\", \"30-13=\", \"6+70=\", \"44-37=\", \"19+36=\"),\n    @(\"9+29=\", \"93-43=\", \"22+41=\", \"65-25=\", \"58-39=\"),\n    @(\"12+43=\", \"99-69=\", \"26-3=\", \"83-50=\", \"35-16=\"),\n    @(\"78-35=\", \"4+41=\", \"25+11=\", \"58+21=\", \"72+18=\"),\n    @(\"84-3=\", \"11+38=\", \"86-72=\", \"65-31=\", \"75-3=\"),\n    @(\"96-76=\", \"22+76=\", \"26+39=\", \"12+5=\", \"68+15=\"),\n    @(\"34-17=\", \"80-77=\", \"0+92=\", \"94-34=\", \"3+18=\"),\n    @(\"76-58=\", \"17+13=\", \"51-50=\", \"9+10=\", \"58+9=\"),\n    @(\"20+63=\", \"34+5=\", \"1-0=\", \"6+75=\", \"64-47=\"),\n    @(\"44-35=\", \"79-12=\", \"97-85=\", \"88-82=\", \"4+48=\"),\n    @(\"20-16=\", \"37+51=\", \"5+71=\", \"94-43=\", \"56-14=\"),\n    @(\"25+70=\", \"9+85=\", \"4+24=\", \"91-47=\", \"92-91=\"),\n    @(\"69-11=\", \"47-4=\", \"42+20=\", \"42+25=\", \"34+22=\")\n)\n\n$table = $d.Tables(1)\nfor ($r = 0; $r -lt $newValues.Length; $r++) {\n    for ($c = 0; $c -lt $newValues[$r].Length; $c++) {\n        $cell = $table.Cell($r + 1, $c + 1)\n        $cell.Range.Text = $newValues[$r][$c]\n    }\n}\n"}
